$wb = $excel.ActiveWorkbook

# --- Sheet "Version": update Start Time / End Time strings ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = '08:49:11 03-21-2017 Tuesday Pacific Daylight Time'
$wsVersion.Range("C11").Value = '08:49:46 03-21-2017 Tuesday Pacific Daylight Time'

# --- Sheet "3_3": update Deep Sleep / BT Idle stats + raw data ---
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = 0.59
$ws33.Range("D3").Value = 0.618
$ws33.Range("E3").Value = 0.5659999999999999
$ws33.Range("F3").Value = 0.016
$ws33.Range("H3").Value = '0.58442,0.60957,0.57188,0.583667,0.609672,0.580704,0.576336,0.602138,0.614027,0.574672,0.588123,0.609387,0.576513,0.575317,0.577355,0.608368,0.595277,0.572437,0.580229,0.612084,0.575596,0.600562,0.617464,0.591187,0.578761,0.573082,0.61811,0.610684,0.577267,0.576248,0.611248,0.615059,0.573361,0.583578,0.580046,0.584135,0.597682,0.577444,0.577552,0.614027,0.606411,0.576798,0.60724,0.586459,0.575678,0.613558,0.601398,0.575032,0.577641,0.608463,0.603347,0.580704,0.584692,0.580426,0.607437,0.608836,0.575507,0.582192,0.600936,0.567702,0.588769,0.611323,0.605576,0.579767,0.615052,0.57671,0.566221,0.606051,0.589801,0.580616,0.580698,0.607342,0.584971,0.576526,0.575956,0.607539,0.587206,0.574665,0.581255,0.608184,0.578938,0.571975,0.610039,0.611805,0.613286,0.569828,0.578748,0.615236,0.584135,0.579665,0.582552,0.579964,0.600474,0.572715,0.584977,0.613748,0.572063,0.578558,0.574101,0.614489'
$ws33.Range("C4").Value = 9.831
$ws33.Range("D4").Value = 9.875
$ws33.Range("E4").Value = 9.808
$ws33.Range("F4").Value = 0.014
$ws33.Range("H4").Value = '9.866544,9.869235,9.875084,9.867557,9.845002,9.854839,9.859017,9.864683,9.855206,9.842496,9.846857,9.856979,9.850661,9.835621,9.829778,9.84435,9.84685,9.834221,9.811578,9.824384,9.838304,9.846762,9.824669,9.814921,9.825322,9.832842,9.828949,9.818175,9.821694,9.824201,9.83757,9.825967,9.817985,9.82096,9.836178,9.8281,9.821232,9.817713,9.824309,9.83507,9.825777,9.822536,9.810091,9.824201,9.832006,9.828284,9.809717,9.826246,9.841558,9.827917,9.819472,9.817067,9.825322,9.832462,9.827917,9.818909,9.816232,9.833671,9.839615,9.822618,9.820953,9.824194,9.838311,9.836551,9.820953,9.81204,9.830607,9.83433,9.832937,9.816871,9.81825,9.836171,9.837387,9.822618,9.813623,9.828847,9.832652,9.83312,9.812699,9.822169,9.829873,9.834418,9.831164,9.81863,9.821137,9.834228,9.834235,9.828562,9.808141,9.82522,9.832475,9.836184,9.818155,9.817061,9.821979,9.838963,9.826531,9.821884,9.823277,9.838026'

# --- Sheet "1_8": update Deep Sleep / BT Idle stats + raw data ---
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = -0.001
$ws18.Range("D3").Value = 0.011
$ws18.Range("E3").Value = -0.011
$ws18.Range("F3").Value = 0.004
$ws18.Range("H3").Value = '-0.001064,-0.004413,-0.000528,-0.001071,0.001704,-0.005523,-0.005245,-0.002754,0.000348,-0.000793,0.002814,-0.001354,-0.001077,0.001143,-0.002187,0.004201,-0.004117,0.000588,-0.001915,0.002814,0.002524,0.000853,0.002826,0.003085,3.3e-05,3.3e-05,-0.007194,-0.002205,0.000613,0.006452,-0.006337,4.5e-05,-0.004931,-0.003007,-0.003549,-0.003556,-0.001077,0.002265,7e-05,-0.002723,-0.003574,-0.004968,-0.001379,0.003653,-0.000799,-0.005794,0.006964,-0.001607,0.000878,-0.006078,-0.002483,-0.004431,-0.002193,-0.011048,-0.003574,-0.001632,-0.000522,-0.00022,0.001704,-0.001348,-0.001064,-0.005523,-0.001058,-0.000824,-0.007712,0.00475,-0.003013,0.006458,-0.000509,0.001686,-0.004135,3.3e-05,-0.003019,0.00118,0.003085,-0.001613,-0.006343,-0.003001,2.1e-05,-0.001626,-0.004931,-0.001089,-0.009944,-0.004394,-0.006633,-0.000232,-0.002181,0.000866,-0.004968,4.5e-05,-0.002193,0.011145,-0.000799,-0.000503,0.003079,0.000878,-0.002199,-0.003303,0.001414,0.007248'
$ws18.Range("C4").Value = -0.001
$ws18.Range("D4").Value = 0.008
$ws18.Range("E4").Value = -0.01
$ws18.Range("F4").Value = 0.004
$ws18.Range("H4").Value = '0.000872,0.000588,-0.003864,0.00118,-0.001058,-0.006071,-0.002717,0.003924,-0.008001,-0.001083,-0.005782,-0.001083,0.000329,-0.004943,-0.005806,-0.000226,-0.006891,-0.002187,0.000298,0.00253,0.002,-0.003007,0.000298,-0.00136,-0.006349,0.00533,0.004195,3.9e-05,-0.0058,0.001988,0.000329,-0.001083,-0.002717,0.001445,-0.010215,-0.003568,0.001162,0.00253,-0.001909,-0.001626,0.001143,0.002518,0.001717,-0.002168,-0.002174,-0.007181,0.001414,-0.00244,-0.003001,0.001988,2.7e-05,-0.000238,0.001125,0.002518,0.005589,0.000329,2.1e-05,-0.000793,0.000866,0.002802,0.000323,0.005854,-0.005794,0.003647,0.000329,0.004756,0.001982,0.00311,-0.005233,-0.003586,0.002259,-0.001354,4.5e-05,-0.003845,-0.001885,-0.001885,-0.005541,-0.004968,-0.000503,-0.002464,0.003653,-0.00633,0.008351,-0.000793,0.000872,-0.003013,-0.005788,-0.006626,0.001717,0.001723,0.000304,5.8e-05,-0.006916,-0.007724,-0.002211,-0.001077,0.002826,0.000619,0.001969,-0.000516'

